# Update countries & provincias Spain
# Refresh COVID-19 case numbers for a batch of countries on the "Pais" sheet
# (columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#  F=Casos criticos, G=Muertes hoy, H=Muertes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Siria's total cases (4883) now exceed Republica de Africa Central's (4854),
# so the two countries swap ranking positions: Siria moves up to row 135 and
# Republica de Africa Central moves down to row 136.
$ws.Range("A135").Value = "Siria"
$ws.Range("A136").Value = "Republica de Africa Central"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 8140287
$ws.Range("C4").Value = 50034
$ws.Range("D4").Value = 5264075
$ws.Range("E4").Value = 2654495
$ws.Range("G4").Value = 844
$ws.Range("H4").Value = 221717

# India (row 5)
$ws.Range("B5").Value = 7305070
$ws.Range("C5").Value = 67988
$ws.Range("D5").Value = 6379428
$ws.Range("E5").Value = 814331
$ws.Range("G5").Value = 694
$ws.Range("H5").Value = 111311

# Brasil (row 6)
$ws.Range("B6").Value = 5140863
$ws.Range("C6").Value = 26040
$ws.Range("D6").Value = 4568813
$ws.Range("E6").Value = 420303
$ws.Range("G6").Value = 684
$ws.Range("H6").Value = 151747

# Sudafrica (row 14)
$ws.Range("B14").Value = 696414
$ws.Range("C14").Value = 1877
$ws.Range("D14").Value = 626898
$ws.Range("E14").Value = 51365
$ws.Range("G14").Value = 123
$ws.Range("H14").Value = 18151

# Alemania (row 23)
$ws.Range("B23").Value = 341742
$ws.Range("C23").Value = 6063
$ws.Range("E23").Value = 50071
$ws.Range("G23").Value = 31
$ws.Range("H23").Value = 9771

# Canada (row 30)
$ws.Range("B30").Value = 188984
$ws.Range("C30").Value = 2103
$ws.Range("D30").Value = 159045
$ws.Range("E30").Value = 20276

# Barein (row 57)
$ws.Range("B57").Value = 76621
$ws.Range("C57").Value = 349
$ws.Range("D57").Value = 72561
$ws.Range("E57").Value = 3773

# Uzbekistan (row 60)
$ws.Range("B60").Value = 61950
$ws.Range("C60").Value = 308
$ws.Range("D60").Value = 58951
$ws.Range("E60").Value = 2485
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 514

# Bulgaria (row 85)
$ws.Range("B85").Value = 26593
$ws.Range("C85").Value = 819
$ws.Range("D85").Value = 16489
$ws.Range("E85").Value = 9175
$ws.Range("G85").Value = 6
$ws.Range("H85").Value = 929

# Gabon (row 114)
$ws.Range("B114").Value = 8869
$ws.Range("C114").Value = 9
$ws.Range("D114").Value = 8395
$ws.Range("E114").Value = 420

# Angola (row 120)
$ws.Range("B120").Value = 6846
$ws.Range("C120").Value = 166
$ws.Range("D120").Value = 2801
$ws.Range("E120").Value = 3818
$ws.Range("G120").Value = 5
$ws.Range("H120").Value = 227

# Republica de Yibuti (row 125)
$ws.Range("B125").Value = 5440
$ws.Range("C125").Value = 12
$ws.Range("D125").Value = 5366
$ws.Range("E125").Value = 13

# Trinidad yTobago (row 131)
$ws.Range("B131").Value = 5154
$ws.Range("C131").Value = 27
$ws.Range("D131").Value = 3452
$ws.Range("E131").Value = 1609

# Guinea Ecuatorial (row 133)
$ws.Range("B133").Value = 5068
$ws.Range("C133").Value = 2
$ws.Range("E133").Value = 31

# Siria (row 135, after the ranking swap above)
$ws.Range("B135").Value = 4883
$ws.Range("C135").Value = 57
$ws.Range("D135").Value = 1389
$ws.Range("E135").Value = 3260
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 234

# Republica de Africa Central (row 136, after the ranking swap above)
$ws.Range("B136").Value = 4854
$ws.Range("D136").Value = 1924
$ws.Range("E136").Value = 2868
$ws.Range("H136").Value = 62

# Aruba (row 138)
$ws.Range("B138").Value = 4255
$ws.Range("C138").Value = 26
$ws.Range("D138").Value = 3875
$ws.Range("E138").Value = 348
$ws.Range("H138").Value = 32

# Sierra Leona (row 157)
$ws.Range("B157").Value = 2315
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 1743
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 73

# Liberia (row 165)
$ws.Range("B165").Value = 1372
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 1254
$ws.Range("E165").Value = 36

# Niger (row 167)
$ws.Range("B167").Value = 1205
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 12
